# Update the "New Orleans_A" team-specific transition-probability matrix with
# freshly computed empirical values (team specific time data).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Af0)
$ws.Range("B2").Value = 0.1778975741239892
$ws.Range("C2").Value = 0.5983827493261455
$ws.Range("P2").Value = 0.1212938005390836
$ws.Range("S2").Value = 0.1024258760107817

# Row 3 (Af1)
$ws.Range("B3").Value = 0.01357466063348416
$ws.Range("C3").Value = 0.01809954751131222
$ws.Range("J3").Value = 0.009049773755656109
$ws.Range("P3").Value = 0.7330316742081447
$ws.Range("S3").Value = 0.2262443438914027

# Row 4 (Af2)
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.6730769230769231
$ws.Range("S4").Value = 0.2884615384615384

# Row 6 (Ai0)
$ws.Range("B6").Value = 0.0728744939271255
$ws.Range("D6").Value = 0.02024291497975709
$ws.Range("E6").Value = 0.008097165991902834
$ws.Range("F6").Value = 0.03643724696356275
$ws.Range("J6").Value = 0.3319838056680162
$ws.Range("O6").Value = 0.03238866396761134
$ws.Range("Q6").Value = 0.1174089068825911
$ws.Range("R6").Value = 0.05263157894736842
$ws.Range("S6").Value = 0.3279352226720648

# Row 7 (Ai1)
$ws.Range("B7").Value = 0.1404255319148936
$ws.Range("D7").Value = 0.01276595744680851
$ws.Range("F7").Value = 0.02978723404255319
$ws.Range("J7").Value = 0.1446808510638298
$ws.Range("O7").Value = 0.03404255319148936
$ws.Range("Q7").Value = 0.1361702127659574
$ws.Range("R7").Value = 0.08085106382978724
$ws.Range("S7").Value = 0.4212765957446808

# Row 8 (Ai2)
$ws.Range("B8").Value = 0.1047794117647059
$ws.Range("D8").Value = 0.02757352941176471
$ws.Range("F8").Value = 0.05882352941176471
$ws.Range("J8").Value = 0.1139705882352941
$ws.Range("O8").Value = 0.02757352941176471
$ws.Range("Q8").Value = 0.1966911764705882
$ws.Range("R8").Value = 0.07169117647058823
$ws.Range("S8").Value = 0.3988970588235294

# Row 9 (Ai3)
$ws.Range("B9").Value = 0.08888888888888889
$ws.Range("D9").Value = 0.02222222222222222
$ws.Range("E9").Value = 0.005555555555555556
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.1555555555555556
$ws.Range("O9").Value = 0.01111111111111111
$ws.Range("Q9").Value = 0.1777777777777778
$ws.Range("R9").Value = 0.04444444444444445
$ws.Range("S9").Value = 0.4444444444444444

# Row 10 (Ar0)
$ws.Range("B10").Value = 0.1316187594553707
$ws.Range("D10").Value = 0.01966717095310136
$ws.Range("E10").Value = 0.0007564296520423601
$ws.Range("F10").Value = 0.06429652042360061
$ws.Range("J10").Value = 0.1270801815431165
$ws.Range("O10").Value = 0.01059001512859304
$ws.Range("Q10").Value = 0.2488653555219364
$ws.Range("R10").Value = 0.05143721633888049
$ws.Range("S10").Value = 0.3456883509833585

# Row 11 (Bf0)
$ws.Range("G11").Value = 0.1218130311614731
$ws.Range("J11").Value = 0.1218130311614731
$ws.Range("K11").Value = 0.1926345609065156
$ws.Range("L11").Value = 0.5410764872521246
$ws.Range("S11").Value = 0.0226628895184136

# Row 12 (Bf1)
$ws.Range("G12").Value = 0.7941176470588235
$ws.Range("J12").Value = 0.1372549019607843
$ws.Range("K12").Value = 0.004901960784313725
$ws.Range("L12").Value = 0.02450980392156863
$ws.Range("S12").Value = 0.0392156862745098

# Row 13 (Bf2)
$ws.Range("F13").Value = 0.01587301587301587
$ws.Range("G13").Value = 0.6031746031746031
$ws.Range("J13").Value = 0.2698412698412698
$ws.Range("S13").Value = 0.1111111111111111

# Row 15 (Bi0)
$ws.Range("F15").Value = 0.0182648401826484
$ws.Range("H15").Value = 0.1917808219178082
$ws.Range("I15").Value = 0.0502283105022831
$ws.Range("J15").Value = 0.273972602739726
$ws.Range("K15").Value = 0.0502283105022831
$ws.Range("M15").Value = 0.0091324200913242
$ws.Range("N15").Value = 0.0091324200913242
$ws.Range("O15").Value = 0.0867579908675799
$ws.Range("S15").Value = 0.3105022831050228

# Row 16 (Bi1)
$ws.Range("F16").Value = 0.0375
$ws.Range("H16").Value = 0.1708333333333333
$ws.Range("I16").Value = 0.075
$ws.Range("J16").Value = 0.3416666666666667
$ws.Range("K16").Value = 0.1375
$ws.Range("M16").Value = 0.0125
$ws.Range("O16").Value = 0.05416666666666667
$ws.Range("S16").Value = 0.1708333333333333

# Row 17 (Bi2)
$ws.Range("F17").Value = 0.04299065420560748
$ws.Range("H17").Value = 0.1850467289719626
$ws.Range("I17").Value = 0.08971962616822429
$ws.Range("J17").Value = 0.3551401869158878
$ws.Range("K17").Value = 0.08971962616822429
$ws.Range("M17").Value = 0.02242990654205607
$ws.Range("N17").Value = 0.001869158878504673
$ws.Range("O17").Value = 0.05981308411214954
$ws.Range("S17").Value = 0.1532710280373832

# Row 18 (Bi3)
$ws.Range("F18").Value = 0.02684563758389262
$ws.Range("H18").Value = 0.2147651006711409
$ws.Range("I18").Value = 0.0738255033557047
$ws.Range("J18").Value = 0.3825503355704698
$ws.Range("K18").Value = 0.08053691275167785
$ws.Range("M18").Value = 0.006711409395973154
$ws.Range("O18").Value = 0.1073825503355705
$ws.Range("S18").Value = 0.1073825503355705

# Row 19 (Br0)
$ws.Range("F19").Value = 0.02725366876310273
$ws.Range("H19").Value = 0.23340321453529
$ws.Range("I19").Value = 0.06359189378057302
$ws.Range("J19").Value = 0.3361285814116003
$ws.Range("K19").Value = 0.1208944793850454
$ws.Range("M19").Value = 0.03144654088050314
$ws.Range("O19").Value = 0.04542278127183787
$ws.Range("S19").Value = 0.1418588399720475
